$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "reads per taxa" counts for samples 2-44 (columns B:H).
$data = New-Object 'object[,]' 43,7
$data[0,0] = 20535
$data[0,1] = 164
$data[0,2] = 0
$data[0,3] = 4
$data[0,4] = 2
$data[0,5] = 0
$data[0,6] = 12683
$data[1,0] = 23389
$data[1,1] = 322
$data[1,2] = 0
$data[1,3] = 2
$data[1,4] = 6
$data[1,5] = 0
$data[1,6] = 9661
$data[2,0] = 24813
$data[2,1] = 140
$data[2,2] = 0
$data[2,3] = 0
$data[2,4] = 0
$data[2,5] = 0
$data[2,6] = 8436
$data[3,0] = 24229
$data[3,1] = 478
$data[3,2] = 0
$data[3,3] = 5
$data[3,4] = 2
$data[3,5] = 0
$data[3,6] = 8666
$data[4,0] = 22092
$data[4,1] = 108
$data[4,2] = 0
$data[4,3] = 3
$data[4,4] = 9
$data[4,5] = 0
$data[4,6] = 11164
$data[5,0] = 24012
$data[5,1] = 30
$data[5,2] = 1
$data[5,3] = 4
$data[5,4] = 32
$data[5,5] = 0
$data[5,6] = 9299
$data[6,0] = 26612
$data[6,1] = 82
$data[6,2] = 1
$data[6,3] = 0
$data[6,4] = 5
$data[6,5] = 0
$data[6,6] = 6678
$data[7,0] = 24606
$data[7,1] = 54
$data[7,2] = 3
$data[7,3] = 3
$data[7,4] = 10
$data[7,5] = 0
$data[7,6] = 8723
$data[8,0] = 25108
$data[8,1] = 42
$data[8,2] = 0
$data[8,3] = 7
$data[8,4] = 4
$data[8,5] = 0
$data[8,6] = 8223
$data[9,0] = 24875
$data[9,1] = 212
$data[9,2] = 0
$data[9,3] = 12
$data[9,4] = 6
$data[9,5] = 0
$data[9,6] = 8290
$data[10,0] = 25964
$data[10,1] = 192
$data[10,2] = 0
$data[10,3] = 0
$data[10,4] = 0
$data[10,5] = 0
$data[10,6] = 7231
$data[11,0] = 22434
$data[11,1] = 252
$data[11,2] = 0
$data[11,3] = 13
$data[11,4] = 5
$data[11,5] = 0
$data[11,6] = 10686
$data[12,0] = 22962
$data[12,1] = 80
$data[12,2] = 0
$data[12,3] = 7
$data[12,4] = 17
$data[12,5] = 0
$data[12,6] = 10315
$data[13,0] = 24938
$data[13,1] = 69
$data[13,2] = 0
$data[13,3] = 5
$data[13,4] = 9
$data[13,5] = 0
$data[13,6] = 8349
$data[14,0] = 21475
$data[14,1] = 139
$data[14,2] = 1
$data[14,3] = 0
$data[14,4] = 3
$data[14,5] = 0
$data[14,6] = 11786
$data[15,0] = 27596
$data[15,1] = 156
$data[15,2] = 0
$data[15,3] = 0
$data[15,4] = 120
$data[15,5] = 0
$data[15,6] = 5518
$data[16,0] = 27038
$data[16,1] = 81
$data[16,2] = 0
$data[16,3] = 2
$data[16,4] = 12
$data[16,5] = 0
$data[16,6] = 6264
$data[17,0] = 22392
$data[17,1] = 246
$data[17,2] = 1
$data[17,3] = 3
$data[17,4] = 5
$data[17,5] = 0
$data[17,6] = 10734
$data[18,0] = 26824
$data[18,1] = 397
$data[18,2] = 10
$data[18,3] = 0
$data[18,4] = 5
$data[18,5] = 0
$data[18,6] = 6154
$data[19,0] = 24065
$data[19,1] = 248
$data[19,2] = 0
$data[19,3] = 4
$data[19,4] = 0
$data[19,5] = 0
$data[19,6] = 9061
$data[20,0] = 25571
$data[20,1] = 327
$data[20,2] = 0
$data[20,3] = 1
$data[20,4] = 15
$data[20,5] = 0
$data[20,6] = 7488
$data[21,0] = 23767
$data[21,1] = 134
$data[21,2] = 0
$data[21,3] = 0
$data[21,4] = 2
$data[21,5] = 0
$data[21,6] = 9500
$data[22,0] = 23112
$data[22,1] = 706
$data[22,2] = 0
$data[22,3] = 0
$data[22,4] = 1
$data[22,5] = 0
$data[22,6] = 9559
$data[23,0] = 25833
$data[23,1] = 191
$data[23,2] = 0
$data[23,3] = 1
$data[23,4] = 1
$data[23,5] = 0
$data[23,6] = 7360
$data[24,0] = 12734
$data[24,1] = 349
$data[24,2] = 0
$data[24,3] = 0
$data[24,4] = 1
$data[24,5] = 0
$data[24,6] = 20296
$data[25,0] = 24638
$data[25,1] = 288
$data[25,2] = 0
$data[25,3] = 5
$data[25,4] = 34
$data[25,5] = 0
$data[25,6] = 8419
$data[26,0] = 24406
$data[26,1] = 91
$data[26,2] = 0
$data[26,3] = 0
$data[26,4] = 14
$data[26,5] = 0
$data[26,6] = 8877
$data[27,0] = 26677
$data[27,1] = 87
$data[27,2] = 0
$data[27,3] = 2
$data[27,4] = 0
$data[27,5] = 0
$data[27,6] = 6622
$data[28,0] = 27174
$data[28,1] = 77
$data[28,2] = 0
$data[28,3] = 2
$data[28,4] = 1
$data[28,5] = 0
$data[28,6] = 6136
$data[29,0] = 21065
$data[29,1] = 334
$data[29,2] = 0
$data[29,3] = 2
$data[29,4] = 49
$data[29,5] = 1
$data[29,6] = 11943
$data[30,0] = 20911
$data[30,1] = 136
$data[30,2] = 1
$data[30,3] = 0
$data[30,4] = 8
$data[30,5] = 1
$data[30,6] = 12334
$data[31,0] = 25908
$data[31,1] = 148
$data[31,2] = 0
$data[31,3] = 5
$data[31,4] = 3
$data[31,5] = 0
$data[31,6] = 7331
$data[32,0] = 25203
$data[32,1] = 326
$data[32,2] = 2
$data[32,3] = 0
$data[32,4] = 5
$data[32,5] = 0
$data[32,6] = 7842
$data[33,0] = 25340
$data[33,1] = 52
$data[33,2] = 0
$data[33,3] = 2
$data[33,4] = 5
$data[33,5] = 0
$data[33,6] = 7994
$data[34,0] = 27702
$data[34,1] = 37
$data[34,2] = 0
$data[34,3] = 0
$data[34,4] = 9
$data[34,5] = 0
$data[34,6] = 5654
$data[35,0] = 27143
$data[35,1] = 111
$data[35,2] = 0
$data[35,3] = 0
$data[35,4] = 5
$data[35,5] = 0
$data[35,6] = 6137
$data[36,0] = 25509
$data[36,1] = 28
$data[36,2] = 1
$data[36,3] = 2
$data[36,4] = 3
$data[36,5] = 0
$data[36,6] = 7859
$data[37,0] = 27590
$data[37,1] = 243
$data[37,2] = 0
$data[37,3] = 1
$data[37,4] = 5
$data[37,5] = 0
$data[37,6] = 5564
$data[38,0] = 27498
$data[38,1] = 140
$data[38,2] = 0
$data[38,3] = 3
$data[38,4] = 4
$data[38,5] = 0
$data[38,6] = 5761
$data[39,0] = 27243
$data[39,1] = 66
$data[39,2] = 0
$data[39,3] = 2
$data[39,4] = 2
$data[39,5] = 0
$data[39,6] = 6077
$data[40,0] = 27666
$data[40,1] = 80
$data[40,2] = 0
$data[40,3] = 1
$data[40,4] = 3
$data[40,5] = 0
$data[40,6] = 5650
$data[41,0] = 25301
$data[41,1] = 639
$data[41,2] = 1
$data[41,3] = 3
$data[41,4] = 0
$data[41,5] = 0
$data[41,6] = 7456
$data[42,0] = 28462
$data[42,1] = 32
$data[42,2] = 0
$data[42,3] = 2
$data[42,4] = 5
$data[42,5] = 0
$data[42,6] = 4895

$ws.Range("B2:H44").Value2 = $data
